# regenerate orders with updated distance/sizes
#
# The experiment's distance conditions and one of the size conditions were
# renamed:
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31
# (S20 and S25 are unchanged.)
#
# These tokens show up embedded inside many of the text values on the sheet
# (Condition, Filename_Left, Filename_Right, Distance, Size columns, etc.),
# so walk every used cell and rewrite any string value that contains one of
# the old tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$startRow = $used.Row()
$startCol = $used.Column()
$rowCount = $used.Rows.Count()
$colCount = $used.Columns.Count()

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
        $val = $cell.Value()

        if ($val -is [string]) {
            $newVal = $val.Replace("D64", "D69").Replace("D51", "D55").Replace("D80", "D86").Replace("S30", "S31")

            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
